$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Вопросы" (sheet1): insert a new column B "Номер_Вопроса" with values
# 1,2,3 before the existing "Вопрос" / "Дата создания" columns, which shift
# right by one column (B->C, C->D).
# ---------------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Вопросы")

# Move existing column C (Дата создания) to D - formats first, then values
$wsQ.Range("C1:C4").Copy()
$wsQ.Range("D1:D4").PasteSpecial(-4122) # xlPasteFormats
$wsQ.Range("C1:C4").Copy()
$wsQ.Range("D1:D4").PasteSpecial(-4163) # xlPasteValues

# Move existing column B (Вопрос) to C - formats first, then values
$wsQ.Range("B1:B4").Copy()
$wsQ.Range("C1:C4").PasteSpecial(-4122) # xlPasteFormats
$wsQ.Range("B1:B4").Copy()
$wsQ.Range("C1:C4").PasteSpecial(-4163) # xlPasteValues

# Build the new column B: header with the same style as the other headers,
# plus a simple incrementing question number per row.
$wsQ.Range("C1").Copy()
$wsQ.Range("B1").PasteSpecial(-4122) # xlPasteFormats
$wsQ.Range("B1").Value = "Номер_Вопроса"

$wsQ.Range("B2").Value = 1
$wsQ.Range("B3").Value = 2
$wsQ.Range("B4").Value = 3

# ---------------------------------------------------------------------------
# Sheet "Ответы" (sheet2): refresh the answer-set statistics - updated
# unique_key (B), normalized Оценка (D) and re-pointed ПВИ (G) values.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Ответы")

$answerRows = @(
    @{ Row = 2;  Key = 179182; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" },
    @{ Row = 3;  Key = 179182; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" },
    @{ Row = 4;  Key = 179182; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" },
    @{ Row = 5;  Key = 212300; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП №2)" },
    @{ Row = 6;  Key = 212300; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП №2)" },
    @{ Row = 7;  Key = 212300; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП №2)" },
    @{ Row = 8;  Key = 13679;  Score = 1; Pvi = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)" },
    @{ Row = 9;  Key = 13679;  Score = 1; Pvi = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)" },
    @{ Row = 10; Key = 13679;  Score = 1; Pvi = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)" },
    @{ Row = 11; Key = 146533; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" },
    @{ Row = 12; Key = 146533; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" },
    @{ Row = 13; Key = 146533; Score = 1; Pvi = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)" }
)

foreach ($r in $answerRows) {
    $wsA.Range("B$($r.Row)").Value = $r.Key
    $wsA.Range("D$($r.Row)").Value = $r.Score
    $wsA.Range("G$($r.Row)").Value = $r.Pvi
}

Write-Output "Workbook updated"
